$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.218.59"
$ws.Range("E2").Value = "  +5.68%  "
$ws.Range("D3").Value = "3.519.82"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.523.67"
$ws.Range("E8").Value = "  +3.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.126"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.13%  "
$ws.Range("D13").Value = "4.128.39"
$ws.Range("E13").Value = "  +3.31%  "
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000181"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.25%  "
$ws.Range("D17").Value = "67.085.37"
$ws.Range("E17").Value = "  +5.32%  "
$ws.Range("D18").Value = "3.530.10"
$ws.Range("E18").Value = "  +3.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "397.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.57%  "
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000128"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.84%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +3.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.24%  "
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.909"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0750"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.09%  "
$ws.Range("E42").Value = "  +5.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("D44").Value = "2.829.45"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "27.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0316"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "353.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.96%  "
